$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Header text updates (rich-text shared strings)
#    A8:  "Volume 31   Number  51"  ->  "...52"   (Volume number bump)
#    C9:  "Report Covering the Week  12/16/2024  Through  12/22/2024"
#         -> "...12/23/2024  Through  12/29/2024"  (week rolled forward)
# ------------------------------------------------------------------

$volCell = $ws.Range("A8")
$volText = $volCell.Value()
$volIdx = $volText.IndexOf("51")
$volChars = $volCell.Characters($volIdx + 1, 2)
$volChars.Text = "52"

$weekCell = $ws.Range("C9")
$weekText = $weekCell.Value()

$firstIdx = $weekText.IndexOf("12/16/2024")
$firstChars = $weekCell.Characters($firstIdx + 1, 10)
$firstChars.Text = "12/23/2024"

# re-read after the first edit so the second date's offset is current
$weekText2 = $weekCell.Value()
$secondIdx = $weekText2.IndexOf("12/22/2024")
$secondChars = $weekCell.Characters($secondIdx + 1, 10)
$secondChars.Text = "12/29/2024"

# ------------------------------------------------------------------
# 2) Row 33 D/E cells switch from the "N/A" / "***.*" text placeholders
#    to real numbers -> copy number format from an already-numeric
#    sibling cell so the style matches (s=14 / s=15), then set value.
# ------------------------------------------------------------------

$ws.Range("D33").NumberFormat = $ws.Range("C15").NumberFormat()
$ws.Range("D33").Value = 2

$ws.Range("E33").NumberFormat = $ws.Range("H14").NumberFormat()
$ws.Range("E33").Value = -100

# ------------------------------------------------------------------
# 3) Weekly crime statistics refresh (rows 14-31, 33)
# ------------------------------------------------------------------

    # Row 14
    $ws.Range("F14").Value = 3
    $ws.Range("G14").Value = 2
    $ws.Range("H14").Value = 50
    $ws.Range("I14").Value = 70
    $ws.Range("K14").Value = 6.060606060606
    $ws.Range("L14").Value = -9.090909090909
    $ws.Range("M14").Value = -46.969696969697
    $ws.Range("N14").Value = -85.169491525423

    # Row 15
    $ws.Range("C15").Value = 3
    $ws.Range("D15").Value = 3
    $ws.Range("E15").Value = 0
    $ws.Range("F15").Value = 19
    $ws.Range("G15").Value = 18
    $ws.Range("H15").Value = 5.555555555555
    $ws.Range("I15").Value = 247
    $ws.Range("J15").Value = 216
    $ws.Range("K15").Value = 14.351851851851
    $ws.Range("L15").Value = -3.13725490196
    $ws.Range("M15").Value = 9.29203539823
    $ws.Range("N15").Value = -58.206429780033

    # Row 16
    $ws.Range("C16").Value = 26
    $ws.Range("D16").Value = 53
    $ws.Range("E16").Value = -50.943396226415
    $ws.Range("F16").Value = 139
    $ws.Range("G16").Value = 200
    $ws.Range("H16").Value = -30.5
    $ws.Range("I16").Value = 2345
    $ws.Range("J16").Value = 2528
    $ws.Range("K16").Value = -7.238924050632
    $ws.Range("L16").Value = -7.018239492466
    $ws.Range("M16").Value = -36.792452830188
    $ws.Range("N16").Value = -86.008353221957

    # Row 17
    $ws.Range("C17").Value = 52
    $ws.Range("D17").Value = 85
    $ws.Range("E17").Value = -38.823529411764
    $ws.Range("F17").Value = 242
    $ws.Range("G17").Value = 301
    $ws.Range("H17").Value = -19.601328903654
    $ws.Range("I17").Value = 4209
    $ws.Range("J17").Value = 4240
    $ws.Range("K17").Value = -0.731132075471
    $ws.Range("L17").Value = 2.558479532163
    $ws.Range("M17").Value = 29.269041769041
    $ws.Range("N17").Value = -49.659131682813

    # Row 18
    $ws.Range("C18").Value = 29
    $ws.Range("D18").Value = 41
    $ws.Range("E18").Value = -29.268292682926
    $ws.Range("G18").Value = 146
    $ws.Range("H18").Value = -26.027397260274
    $ws.Range("I18").Value = 1872
    $ws.Range("J18").Value = 2033
    $ws.Range("K18").Value = -7.919331037875
    $ws.Range("L18").Value = -20.87912087912
    $ws.Range("M18").Value = -41.917468197331
    $ws.Range("N18").Value = -84.459571642038

    # Row 19
    $ws.Range("C19").Value = 56
    $ws.Range("D19").Value = 90
    $ws.Range("E19").Value = -37.777777777777
    $ws.Range("F19").Value = 284
    $ws.Range("G19").Value = 419
    $ws.Range("H19").Value = -32.219570405727
    $ws.Range("I19").Value = 5256
    $ws.Range("J19").Value = 5800
    $ws.Range("K19").Value = -9.379310344827
    $ws.Range("L19").Value = -11.380880121396
    $ws.Range("M19").Value = 21.385681293302
    $ws.Range("N19").Value = -23.770848440899

    # Row 20
    $ws.Range("C20").Value = 18
    $ws.Range("D20").Value = 29
    $ws.Range("E20").Value = -37.931034482758
    $ws.Range("F20").Value = 83
    $ws.Range("G20").Value = 141
    $ws.Range("H20").Value = -41.134751773049
    $ws.Range("I20").Value = 1704
    $ws.Range("J20").Value = 1865
    $ws.Range("K20").Value = -8.632707774798
    $ws.Range("L20").Value = -8.189655172413
    $ws.Range("M20").Value = 19.662921348314
    $ws.Range("N20").Value = -82.138364779874

    # Row 21
    $ws.Range("C21").Value = 185
    $ws.Range("D21").Value = 301
    $ws.Range("E21").Value = -38.538205980066
    $ws.Range("F21").Value = 878
    $ws.Range("G21").Value = 1227
    $ws.Range("H21").Value = -28.443357783211
    $ws.Range("I21").Value = 15703
    $ws.Range("J21").Value = 16748
    $ws.Range("K21").Value = -6.239550991163
    $ws.Range("L21").Value = -8.2286248612
    $ws.Range("M21").Value = -3.668486595914
    $ws.Range("N21").Value = -71.274124211104

    # Row 22
    $ws.Range("C22").Value = 4
    $ws.Range("D22").Value = 9
    $ws.Range("E22").Value = -55.555555555555
    $ws.Range("G22").Value = 28
    $ws.Range("H22").Value = -25
    $ws.Range("I22").Value = 291
    $ws.Range("J22").Value = 309
    $ws.Range("K22").Value = -5.825242718446
    $ws.Range("L22").Value = -18.258426966292
    $ws.Range("M22").Value = -34.89932885906

    # Row 23
    $ws.Range("C23").Value = 19
    $ws.Range("D23").Value = 26
    $ws.Range("E23").Value = -26.923076923076
    $ws.Range("F23").Value = 90
    $ws.Range("G23").Value = 99
    $ws.Range("H23").Value = -9.090909090909
    $ws.Range("I23").Value = 1398
    $ws.Range("J23").Value = 1554
    $ws.Range("K23").Value = -10.03861003861
    $ws.Range("L23").Value = -6.862091938707
    $ws.Range("M23").Value = 19.181585677749

    # Row 24
    $ws.Range("C24").Value = 206
    $ws.Range("D24").Value = 188
    $ws.Range("E24").Value = 9.574468085106
    $ws.Range("F24").Value = 999
    $ws.Range("G24").Value = 925
    $ws.Range("H24").Value = 8
    $ws.Range("I24").Value = 12480
    $ws.Range("J24").Value = 12176
    $ws.Range("K24").Value = 2.496714848883
    $ws.Range("L24").Value = -6.214774179003
    $ws.Range("M24").Value = 19.460132095338

    # Row 25
    $ws.Range("C25").Value = 77
    $ws.Range("D25").Value = 69
    $ws.Range("E25").Value = 11.59420289855
    $ws.Range("F25").Value = 324
    $ws.Range("G25").Value = 357
    $ws.Range("H25").Value = -9.243697478991
    $ws.Range("I25").Value = 5417
    $ws.Range("J25").Value = 4839
    $ws.Range("K25").Value = 11.944616656334
    $ws.Range("L25").Value = -1.937002172338

    # Row 26
    $ws.Range("C26").Value = 104
    $ws.Range("D26").Value = 106
    $ws.Range("E26").Value = -1.88679245283
    $ws.Range("F26").Value = 429
    $ws.Range("G26").Value = 417
    $ws.Range("H26").Value = 2.877697841726
    $ws.Range("I26").Value = 6350
    $ws.Range("J26").Value = 6050
    $ws.Range("K26").Value = 4.95867768595
    $ws.Range("L26").Value = 8.528456674072
    $ws.Range("M26").Value = -18.306959989708

    # Row 27
    $ws.Range("C27").Value = 3
    $ws.Range("D27").Value = 4
    $ws.Range("E27").Value = -25
    $ws.Range("F27").Value = 20
    $ws.Range("G27").Value = 25
    $ws.Range("H27").Value = -20
    $ws.Range("I27").Value = 335
    $ws.Range("J27").Value = 334
    $ws.Range("K27").Value = 0.299401197604
    $ws.Range("L27").Value = -11.842105263157

    # Row 28
    $ws.Range("C28").Value = 12
    $ws.Range("D28").Value = 11
    $ws.Range("E28").Value = 9.090909090909
    $ws.Range("F28").Value = 47
    $ws.Range("G28").Value = 48
    $ws.Range("H28").Value = -2.083333333333
    $ws.Range("I28").Value = 648
    $ws.Range("J28").Value = 618
    $ws.Range("K28").Value = 4.854368932038
    $ws.Range("L28").Value = 7.820299500831

    # Row 29
    $ws.Range("C29").Value = 2
    $ws.Range("D29").Value = 6
    $ws.Range("G29").Value = 17
    $ws.Range("H29").Value = -17.647058823529
    $ws.Range("I29").Value = 237
    $ws.Range("J29").Value = 244
    $ws.Range("K29").Value = -2.868852459016
    $ws.Range("L29").Value = -30.088495575221
    $ws.Range("M29").Value = -52.50501002004
    $ws.Range("N29").Value = -87.161430119176

    # Row 30
    $ws.Range("D30").Value = 4
    $ws.Range("E30").Value = -75
    $ws.Range("F30").Value = 11
    $ws.Range("G30").Value = 14
    $ws.Range("H30").Value = -21.428571428571
    $ws.Range("I30").Value = 195
    $ws.Range("J30").Value = 206
    $ws.Range("K30").Value = -5.339805825242
    $ws.Range("L30").Value = -32.055749128919
    $ws.Range("M30").Value = -51.970443349753
    $ws.Range("N30").Value = -88.260084286574

    # Row 31
    $ws.Range("D31").Value = 2
    $ws.Range("G31").Value = 10
    $ws.Range("H31").Value = -50
    $ws.Range("J31").Value = 92
    $ws.Range("K31").Value = -4.347826086956

    # Row 33
    $ws.Range("G33").Value = 3
    $ws.Range("J33").Value = 24
    $ws.Range("K33").Value = 4.166666666666
    $ws.Range("L33").Value = -16.666666666666
